$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("gasline_ht")
$ws.Cells.Clear()

# ---- Column widths (B:E and F) ----
$ws.Range("B:E").ColumnWidth = 14
$ws.Range("F:F").ColumnWidth = 14

# ---- Header row (block 1: A-F) ----
$ws.Range("A1").Value = "location"
$ws.Range("B1").Value = "spp"
$ws.Range("B1").HorizontalAlignment = -4108
$ws.Range("C1").Value = "species"
$ws.Range("C1").HorizontalAlignment = -4108
$ws.Range("D1").Value = "ht.2019"
$ws.Range("D1").HorizontalAlignment = -4108
$ws.Range("E1").Value = "inc.2019"
$ws.Range("E1").HorizontalAlignment = -4108
$ws.Range("F1").Value = "n"
$ws.Range("F1").HorizontalAlignment = -4108
$ws.Range("F1").NumberFormat = "0"

# ---- Header row (block 2: H-K) ----
$ws.Range("H1").Value = "location"
$ws.Range("I1").Value = "species"
$ws.Range("I1").HorizontalAlignment = -4108
$ws.Range("J1").Value = "ht.2019"
$ws.Range("J1").HorizontalAlignment = -4108
$ws.Range("K1").Value = "n"
$ws.Range("K1").HorizontalAlignment = -4108
$ws.Range("K1").NumberFormat = "0"

# ---- Header row (block 3: M-P) ----
$ws.Range("M1").Value = "location"
$ws.Range("N1").Value = "species"
$ws.Range("N1").HorizontalAlignment = -4108
$ws.Range("O1").Value = "ht.2019"
$ws.Range("O1").HorizontalAlignment = -4108
$ws.Range("P1").Value = "n"
$ws.Range("P1").HorizontalAlignment = -4108
$ws.Range("P1").NumberFormat = "0"

# ---- Row 2: control 316 RM ----
$ws.Range("A2").Value = "control"
$ws.Range("B2").Value = 316
$ws.Range("B2").HorizontalAlignment = -4108
$ws.Range("C2").Value = "RM"
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("D2").Value = 46.6666666666667
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").NumberFormat = "0"
$ws.Range("E2").Value = 15.6666666666667
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").NumberFormat = "0"
$ws.Range("F2").Value = 15
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").NumberFormat = "0"
$ws.Range("H2").Value = "control"
$ws.Range("I2").Value = "RM"
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("J2").Value = 46.6666666666667
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").NumberFormat = "0"
$ws.Range("K2").Value = 15
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").NumberFormat = "0"
$ws.Range("M2").Value = "control"
$ws.Range("N2").Value = "RM"
$ws.Range("N2").HorizontalAlignment = -4108
$ws.Range("O2").Value = 46.6666666666667
$ws.Range("O2").HorizontalAlignment = -4108
$ws.Range("O2").NumberFormat = "0"
$ws.Range("P2").Value = 15
$ws.Range("P2").HorizontalAlignment = -4108
$ws.Range("P2").NumberFormat = "0"
$ws.Range("M2").Interior.Color = 65535

# ---- Row 3: control 356 SV ----
$ws.Range("A3").Value = "control"
$ws.Range("B3").Value = 356
$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("C3").Value = "SV"
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("D3").Value = 32
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").NumberFormat = "0"
$ws.Range("E3").Value = 5
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").NumberFormat = "0"
$ws.Range("F3").Value = 1
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").NumberFormat = "0"
$ws.Range("H3").Value = "control"
$ws.Range("I3").Value = "SV"
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("J3").Value = 32
$ws.Range("J3").HorizontalAlignment = -4108
$ws.Range("J3").NumberFormat = "0"
$ws.Range("K3").Value = 1
$ws.Range("K3").HorizontalAlignment = -4108
$ws.Range("K3").NumberFormat = "0"
$ws.Range("M3").Value = "control"
$ws.Range("N3").Value = "SV"
$ws.Range("N3").HorizontalAlignment = -4108
$ws.Range("O3").Value = 32
$ws.Range("O3").HorizontalAlignment = -4108
$ws.Range("O3").NumberFormat = "0"
$ws.Range("P3").Value = 1
$ws.Range("P3").HorizontalAlignment = -4108
$ws.Range("P3").NumberFormat = "0"
$ws.Range("M3").Interior.Color = 65535

# ---- Row 4: control 372 SB ----
$ws.Range("A4").Value = "control"
$ws.Range("B4").Value = 372
$ws.Range("B4").HorizontalAlignment = -4108
$ws.Range("C4").Value = "SB"
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("D4").Value = 38.818181818181799
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").NumberFormat = "0"
$ws.Range("E4").Value = 19.386363636363601
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("E4").NumberFormat = "0"
$ws.Range("F4").Value = 22
$ws.Range("F4").HorizontalAlignment = -4108
$ws.Range("F4").NumberFormat = "0"
$ws.Range("H4").Value = "control"
$ws.Range("I4").Value = "SB"
$ws.Range("I4").HorizontalAlignment = -4108
$ws.Range("J4").Value = 38.818181818181799
$ws.Range("J4").HorizontalAlignment = -4108
$ws.Range("J4").NumberFormat = "0"
$ws.Range("K4").Value = 22
$ws.Range("K4").HorizontalAlignment = -4108
$ws.Range("K4").NumberFormat = "0"
$ws.Range("M4").Value = "control"
$ws.Range("N4").Value = "SB"
$ws.Range("N4").HorizontalAlignment = -4108
$ws.Range("O4").Value = 38.818181818181799
$ws.Range("O4").HorizontalAlignment = -4108
$ws.Range("O4").NumberFormat = "0"
$ws.Range("P4").Value = 22
$ws.Range("P4").HorizontalAlignment = -4108
$ws.Range("P4").NumberFormat = "0"
$ws.Range("M4").Interior.Color = 65535

# ---- Row 5: control 531 AB ----
$ws.Range("A5").Value = "control"
$ws.Range("B5").Value = 531
$ws.Range("B5").HorizontalAlignment = -4108
$ws.Range("C5").Value = "AB"
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Range("D5").Value = 29.3333333333333
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").NumberFormat = "0"
$ws.Range("E5").Value = 14.2424242424242
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E5").NumberFormat = "0"
$ws.Range("F5").Value = 33
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("F5").NumberFormat = "0"
$ws.Range("H5").Value = "control"
$ws.Range("I5").Value = "AB"
$ws.Range("I5").HorizontalAlignment = -4108
$ws.Range("J5").Value = 29.3333333333333
$ws.Range("J5").HorizontalAlignment = -4108
$ws.Range("J5").NumberFormat = "0"
$ws.Range("K5").Value = 33
$ws.Range("K5").HorizontalAlignment = -4108
$ws.Range("K5").NumberFormat = "0"
$ws.Range("M5").Value = "control"
$ws.Range("N5").Value = "AB"
$ws.Range("N5").HorizontalAlignment = -4108
$ws.Range("O5").Value = 29.3333333333333
$ws.Range("O5").HorizontalAlignment = -4108
$ws.Range("O5").NumberFormat = "0"
$ws.Range("P5").Value = 33
$ws.Range("P5").HorizontalAlignment = -4108
$ws.Range("P5").NumberFormat = "0"
$ws.Range("M5").Interior.Color = 65535

# ---- Row 6: control 541 WA ----
$ws.Range("A6").Value = "control"
$ws.Range("B6").Value = 541
$ws.Range("B6").HorizontalAlignment = -4108
$ws.Range("C6").Value = "WA"
$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("D6").Value = 28.375
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").NumberFormat = "0"
$ws.Range("E6").Value = 9.3125
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").NumberFormat = "0"
$ws.Range("F6").Value = 8
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("F6").NumberFormat = "0"
$ws.Range("H6").Value = "control"
$ws.Range("I6").Value = "WA"
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("J6").Value = 28.375
$ws.Range("J6").HorizontalAlignment = -4108
$ws.Range("J6").NumberFormat = "0"
$ws.Range("K6").Value = 8
$ws.Range("K6").HorizontalAlignment = -4108
$ws.Range("K6").NumberFormat = "0"
$ws.Range("M6").Value = "control"
$ws.Range("N6").Value = "WA"
$ws.Range("N6").HorizontalAlignment = -4108
$ws.Range("O6").Value = 28.375
$ws.Range("O6").HorizontalAlignment = -4108
$ws.Range("O6").NumberFormat = "0"
$ws.Range("P6").Value = 8
$ws.Range("P6").HorizontalAlignment = -4108
$ws.Range("P6").NumberFormat = "0"
$ws.Range("M6").Interior.Color = 65535

# ---- Row 7: control 762 BC ----
$ws.Range("A7").Value = "control"
$ws.Range("B7").Value = 762
$ws.Range("B7").HorizontalAlignment = -4108
$ws.Range("C7").Value = "BC"
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("D7").Value = 33.6666666666667
$ws.Range("D7").HorizontalAlignment = -4108
$ws.Range("D7").NumberFormat = "0"
$ws.Range("E7").Value = 12.3333333333333
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").NumberFormat = "0"
$ws.Range("F7").Value = 3
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F7").NumberFormat = "0"
$ws.Range("H7").Value = "control"
$ws.Range("I7").Value = "BC"
$ws.Range("I7").HorizontalAlignment = -4108
$ws.Range("J7").Value = 33.6666666666667
$ws.Range("J7").HorizontalAlignment = -4108
$ws.Range("J7").NumberFormat = "0"
$ws.Range("K7").Value = 3
$ws.Range("K7").HorizontalAlignment = -4108
$ws.Range("K7").NumberFormat = "0"
$ws.Range("M7").Value = "control"
$ws.Range("N7").Value = "BC"
$ws.Range("N7").HorizontalAlignment = -4108
$ws.Range("O7").Value = 33.6666666666667
$ws.Range("O7").HorizontalAlignment = -4108
$ws.Range("O7").NumberFormat = "0"
$ws.Range("P7").Value = 3
$ws.Range("P7").HorizontalAlignment = -4108
$ws.Range("P7").NumberFormat = "0"
$ws.Range("M7").Interior.Color = 65535

# ---- Row 8: control 833 RO ----
$ws.Range("A8").Value = "control"
$ws.Range("B8").Value = 833
$ws.Range("B8").HorizontalAlignment = -4108
$ws.Range("C8").Value = "RO"
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("D8").Value = 23.772727272727298
$ws.Range("D8").HorizontalAlignment = -4108
$ws.Range("D8").NumberFormat = "0"
$ws.Range("E8").Value = 8.0454545454545396
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").NumberFormat = "0"
$ws.Range("F8").Value = 11
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F8").NumberFormat = "0"
$ws.Range("H8").Value = "control"
$ws.Range("I8").Value = "RO"
$ws.Range("I8").HorizontalAlignment = -4108
$ws.Range("J8").Value = 23.772727272727298
$ws.Range("J8").HorizontalAlignment = -4108
$ws.Range("J8").NumberFormat = "0"
$ws.Range("K8").Value = 11
$ws.Range("K8").HorizontalAlignment = -4108
$ws.Range("K8").NumberFormat = "0"
$ws.Range("M8").Value = "control"
$ws.Range("N8").Value = "RO"
$ws.Range("N8").HorizontalAlignment = -4108
$ws.Range("O8").Value = 23.772727272727298
$ws.Range("O8").HorizontalAlignment = -4108
$ws.Range("O8").NumberFormat = "0"
$ws.Range("P8").Value = 11
$ws.Range("P8").HorizontalAlignment = -4108
$ws.Range("P8").NumberFormat = "0"
$ws.Range("M8").Interior.Color = 65535

# ---- Row 9: interior 316 RM ----
$ws.Range("A9").Value = "interior"
$ws.Range("B9").Value = 316
$ws.Range("B9").HorizontalAlignment = -4108
$ws.Range("C9").Value = "RM"
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("D9").Value = 42.9
$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").NumberFormat = "0"
$ws.Range("E9").Value = 13.1
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").NumberFormat = "0"
$ws.Range("F9").Value = 20
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("F9").NumberFormat = "0"
$ws.Range("H9").Value = "interior"
$ws.Range("I9").Value = "RM"
$ws.Range("I9").HorizontalAlignment = -4108
$ws.Range("J9").Value = 42.9
$ws.Range("J9").HorizontalAlignment = -4108
$ws.Range("J9").NumberFormat = "0"
$ws.Range("K9").Value = 20
$ws.Range("K9").HorizontalAlignment = -4108
$ws.Range("K9").NumberFormat = "0"
$ws.Range("M9").Value = "interior"
$ws.Range("N9").Value = "RM"
$ws.Range("N9").HorizontalAlignment = -4108
$ws.Range("O9").Value = 42.9
$ws.Range("O9").HorizontalAlignment = -4108
$ws.Range("O9").NumberFormat = "0"
$ws.Range("P9").Value = 20
$ws.Range("P9").HorizontalAlignment = -4108
$ws.Range("P9").NumberFormat = "0"

# ---- Row 10: interior 318 SM ----
$ws.Range("A10").Value = "interior"
$ws.Range("B10").Value = 318
$ws.Range("B10").HorizontalAlignment = -4108
$ws.Range("C10").Value = "SM"
$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("D10").Value = 57.15
$ws.Range("D10").HorizontalAlignment = -4108
$ws.Range("D10").NumberFormat = "0"
$ws.Range("E10").Value = 11.26
$ws.Range("E10").HorizontalAlignment = -4108
$ws.Range("E10").NumberFormat = "0"
$ws.Range("F10").Value = 50
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F10").NumberFormat = "0"
$ws.Range("H10").Value = "interior"
$ws.Range("I10").Value = "SM"
$ws.Range("I10").HorizontalAlignment = -4108
$ws.Range("J10").Value = 57.15
$ws.Range("J10").HorizontalAlignment = -4108
$ws.Range("J10").NumberFormat = "0"
$ws.Range("K10").Value = 50
$ws.Range("K10").HorizontalAlignment = -4108
$ws.Range("K10").NumberFormat = "0"
$ws.Range("M10").Value = "interior"
$ws.Range("N10").Value = "SM"
$ws.Range("N10").HorizontalAlignment = -4108
$ws.Range("O10").Value = 57.15
$ws.Range("O10").HorizontalAlignment = -4108
$ws.Range("O10").NumberFormat = "0"
$ws.Range("P10").Value = 50
$ws.Range("P10").HorizontalAlignment = -4108
$ws.Range("P10").NumberFormat = "0"

# ---- Row 11: interior 356 SV ----
$ws.Range("A11").Value = "interior"
$ws.Range("B11").Value = 356
$ws.Range("B11").HorizontalAlignment = -4108
$ws.Range("C11").Value = "SV"
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("D11").Value = 89
$ws.Range("D11").HorizontalAlignment = -4108
$ws.Range("D11").NumberFormat = "0"
$ws.Range("E11").Value = 35.25
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("E11").NumberFormat = "0"
$ws.Range("F11").Value = 2
$ws.Range("F11").HorizontalAlignment = -4108
$ws.Range("F11").NumberFormat = "0"
$ws.Range("H11").Value = "interior"
$ws.Range("I11").Value = "SV"
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("J11").Value = 89
$ws.Range("J11").HorizontalAlignment = -4108
$ws.Range("J11").NumberFormat = "0"
$ws.Range("K11").Value = 2
$ws.Range("K11").HorizontalAlignment = -4108
$ws.Range("K11").NumberFormat = "0"
$ws.Range("M11").Value = "interior"
$ws.Range("N11").Value = "SV"
$ws.Range("N11").HorizontalAlignment = -4108
$ws.Range("O11").Value = 89
$ws.Range("O11").HorizontalAlignment = -4108
$ws.Range("O11").NumberFormat = "0"
$ws.Range("P11").Value = 2
$ws.Range("P11").HorizontalAlignment = -4108
$ws.Range("P11").NumberFormat = "0"

# ---- Row 12: interior 372 SB ----
$ws.Range("A12").Value = "interior"
$ws.Range("B12").Value = 372
$ws.Range("B12").HorizontalAlignment = -4108
$ws.Range("C12").Value = "SB"
$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("D12").Value = 73.5
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("D12").NumberFormat = "0"
$ws.Range("E12").Value = 45
$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E12").NumberFormat = "0"
$ws.Range("F12").Value = 1
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").NumberFormat = "0"
$ws.Range("H12").Value = "interior"
$ws.Range("I12").Value = "SB"
$ws.Range("I12").HorizontalAlignment = -4108
$ws.Range("J12").Value = 73.5
$ws.Range("J12").HorizontalAlignment = -4108
$ws.Range("J12").NumberFormat = "0"
$ws.Range("K12").Value = 1
$ws.Range("K12").HorizontalAlignment = -4108
$ws.Range("K12").NumberFormat = "0"
$ws.Range("M12").Value = "interior"
$ws.Range("N12").Value = "SB"
$ws.Range("N12").HorizontalAlignment = -4108
$ws.Range("O12").Value = 73.5
$ws.Range("O12").HorizontalAlignment = -4108
$ws.Range("O12").NumberFormat = "0"
$ws.Range("P12").Value = 1
$ws.Range("P12").HorizontalAlignment = -4108
$ws.Range("P12").NumberFormat = "0"

# ---- Row 13: interior 531 AB ----
$ws.Range("A13").Value = "interior"
$ws.Range("B13").Value = 531
$ws.Range("B13").HorizontalAlignment = -4108
$ws.Range("C13").Value = "AB"
$ws.Range("C13").HorizontalAlignment = -4108
$ws.Range("D13").Value = 41.909090909090899
$ws.Range("D13").HorizontalAlignment = -4108
$ws.Range("D13").NumberFormat = "0"
$ws.Range("E13").Value = 16.602272727272702
$ws.Range("E13").HorizontalAlignment = -4108
$ws.Range("E13").NumberFormat = "0"
$ws.Range("F13").Value = 44
$ws.Range("F13").HorizontalAlignment = -4108
$ws.Range("F13").NumberFormat = "0"
$ws.Range("H13").Value = "interior"
$ws.Range("I13").Value = "AB"
$ws.Range("I13").HorizontalAlignment = -4108
$ws.Range("J13").Value = 41.909090909090899
$ws.Range("J13").HorizontalAlignment = -4108
$ws.Range("J13").NumberFormat = "0"
$ws.Range("K13").Value = 44
$ws.Range("K13").HorizontalAlignment = -4108
$ws.Range("K13").NumberFormat = "0"
$ws.Range("M13").Value = "interior"
$ws.Range("N13").Value = "AB"
$ws.Range("N13").HorizontalAlignment = -4108
$ws.Range("O13").Value = 41.909090909090899
$ws.Range("O13").HorizontalAlignment = -4108
$ws.Range("O13").NumberFormat = "0"
$ws.Range("P13").Value = 44
$ws.Range("P13").HorizontalAlignment = -4108
$ws.Range("P13").NumberFormat = "0"

# ---- Row 14: interior 541 WA ----
$ws.Range("A14").Value = "interior"
$ws.Range("B14").Value = 541
$ws.Range("B14").HorizontalAlignment = -4108
$ws.Range("C14").Value = "WA"
$ws.Range("C14").HorizontalAlignment = -4108
$ws.Range("D14").Value = 54.682692307692299
$ws.Range("D14").HorizontalAlignment = -4108
$ws.Range("D14").NumberFormat = "0"
$ws.Range("E14").Value = 20.663461538461501
$ws.Range("E14").HorizontalAlignment = -4108
$ws.Range("E14").NumberFormat = "0"
$ws.Range("F14").Value = 52
$ws.Range("F14").HorizontalAlignment = -4108
$ws.Range("F14").NumberFormat = "0"
$ws.Range("H14").Value = "interior"
$ws.Range("I14").Value = "WA"
$ws.Range("I14").HorizontalAlignment = -4108
$ws.Range("J14").Value = 54.682692307692299
$ws.Range("J14").HorizontalAlignment = -4108
$ws.Range("J14").NumberFormat = "0"
$ws.Range("K14").Value = 52
$ws.Range("K14").HorizontalAlignment = -4108
$ws.Range("K14").NumberFormat = "0"
$ws.Range("M14").Value = "interior"
$ws.Range("N14").Value = "WA"
$ws.Range("N14").HorizontalAlignment = -4108
$ws.Range("O14").Value = 54.682692307692299
$ws.Range("O14").HorizontalAlignment = -4108
$ws.Range("O14").NumberFormat = "0"
$ws.Range("P14").Value = 52
$ws.Range("P14").HorizontalAlignment = -4108
$ws.Range("P14").NumberFormat = "0"

# ---- Row 15: interior 621 YP ----
$ws.Range("A15").Value = "interior"
$ws.Range("B15").Value = 621
$ws.Range("B15").HorizontalAlignment = -4108
$ws.Range("C15").Value = "YP"
$ws.Range("C15").HorizontalAlignment = -4108
$ws.Range("D15").Value = 34.5
$ws.Range("D15").HorizontalAlignment = -4108
$ws.Range("D15").NumberFormat = "0"
$ws.Range("E15").Value = 17.5
$ws.Range("E15").HorizontalAlignment = -4108
$ws.Range("E15").NumberFormat = "0"
$ws.Range("F15").Value = 3
$ws.Range("F15").HorizontalAlignment = -4108
$ws.Range("F15").NumberFormat = "0"
$ws.Range("H15").Value = "interior"
$ws.Range("I15").Value = "YP"
$ws.Range("I15").HorizontalAlignment = -4108
$ws.Range("J15").Value = 34.5
$ws.Range("J15").HorizontalAlignment = -4108
$ws.Range("J15").NumberFormat = "0"
$ws.Range("K15").Value = 3
$ws.Range("K15").HorizontalAlignment = -4108
$ws.Range("K15").NumberFormat = "0"
$ws.Range("M15").Value = "interior"
$ws.Range("N15").Value = "YP"
$ws.Range("N15").HorizontalAlignment = -4108
$ws.Range("O15").Value = 34.5
$ws.Range("O15").HorizontalAlignment = -4108
$ws.Range("O15").NumberFormat = "0"
$ws.Range("P15").Value = 3
$ws.Range("P15").HorizontalAlignment = -4108
$ws.Range("P15").NumberFormat = "0"

# ---- Row 16: interior 762 BC ----
$ws.Range("A16").Value = "interior"
$ws.Range("B16").Value = 762
$ws.Range("B16").HorizontalAlignment = -4108
$ws.Range("C16").Value = "BC"
$ws.Range("C16").HorizontalAlignment = -4108
$ws.Range("D16").Value = 80.9444444444444
$ws.Range("D16").HorizontalAlignment = -4108
$ws.Range("D16").NumberFormat = "0"
$ws.Range("E16").Value = 25.5555555555556
$ws.Range("E16").HorizontalAlignment = -4108
$ws.Range("E16").NumberFormat = "0"
$ws.Range("F16").Value = 9
$ws.Range("F16").HorizontalAlignment = -4108
$ws.Range("F16").NumberFormat = "0"
$ws.Range("H16").Value = "interior"
$ws.Range("I16").Value = "BC"
$ws.Range("I16").HorizontalAlignment = -4108
$ws.Range("J16").Value = 80.9444444444444
$ws.Range("J16").HorizontalAlignment = -4108
$ws.Range("J16").NumberFormat = "0"
$ws.Range("K16").Value = 9
$ws.Range("K16").HorizontalAlignment = -4108
$ws.Range("K16").NumberFormat = "0"
$ws.Range("M16").Value = "interior"
$ws.Range("N16").Value = "BC"
$ws.Range("N16").HorizontalAlignment = -4108
$ws.Range("O16").Value = 80.9444444444444
$ws.Range("O16").HorizontalAlignment = -4108
$ws.Range("O16").NumberFormat = "0"
$ws.Range("P16").Value = 9
$ws.Range("P16").HorizontalAlignment = -4108
$ws.Range("P16").NumberFormat = "0"

# ---- Row 17: interior 833 RO ----
$ws.Range("A17").Value = "interior"
$ws.Range("B17").Value = 833
$ws.Range("B17").HorizontalAlignment = -4108
$ws.Range("C17").Value = "RO"
$ws.Range("C17").HorizontalAlignment = -4108
$ws.Range("D17").Value = 42.709677419354797
$ws.Range("D17").HorizontalAlignment = -4108
$ws.Range("D17").NumberFormat = "0"
$ws.Range("E17").Value = 13.1290322580645
$ws.Range("E17").HorizontalAlignment = -4108
$ws.Range("E17").NumberFormat = "0"
$ws.Range("F17").Value = 31
$ws.Range("F17").HorizontalAlignment = -4108
$ws.Range("F17").NumberFormat = "0"
$ws.Range("H17").Value = "interior"
$ws.Range("I17").Value = "RO"
$ws.Range("I17").HorizontalAlignment = -4108
$ws.Range("J17").Value = 42.709677419354797
$ws.Range("J17").HorizontalAlignment = -4108
$ws.Range("J17").NumberFormat = "0"
$ws.Range("K17").Value = 31
$ws.Range("K17").HorizontalAlignment = -4108
$ws.Range("K17").NumberFormat = "0"
$ws.Range("M17").Value = "interior"
$ws.Range("N17").Value = "RO"
$ws.Range("N17").HorizontalAlignment = -4108
$ws.Range("O17").Value = 42.709677419354797
$ws.Range("O17").HorizontalAlignment = -4108
$ws.Range("O17").NumberFormat = "0"
$ws.Range("P17").Value = 31
$ws.Range("P17").HorizontalAlignment = -4108
$ws.Range("P17").NumberFormat = "0"

# ---- Row 18: interior 951 BA ----
$ws.Range("A18").Value = "interior"
$ws.Range("B18").Value = 951
$ws.Range("B18").HorizontalAlignment = -4108
$ws.Range("C18").Value = "BA"
$ws.Range("C18").HorizontalAlignment = -4108
$ws.Range("D18").Value = 59.9375
$ws.Range("D18").HorizontalAlignment = -4108
$ws.Range("D18").NumberFormat = "0"
$ws.Range("E18").Value = 22
$ws.Range("E18").HorizontalAlignment = -4108
$ws.Range("E18").NumberFormat = "0"
$ws.Range("F18").Value = 8
$ws.Range("F18").HorizontalAlignment = -4108
$ws.Range("F18").NumberFormat = "0"
$ws.Range("H18").Value = "interior"
$ws.Range("I18").Value = "BA"
$ws.Range("I18").HorizontalAlignment = -4108
$ws.Range("J18").Value = 59.9375
$ws.Range("J18").HorizontalAlignment = -4108
$ws.Range("J18").NumberFormat = "0"
$ws.Range("K18").Value = 8
$ws.Range("K18").HorizontalAlignment = -4108
$ws.Range("K18").NumberFormat = "0"
$ws.Range("M18").Value = "interior"
$ws.Range("N18").Value = "BA"
$ws.Range("N18").HorizontalAlignment = -4108
$ws.Range("O18").Value = 59.9375
$ws.Range("O18").HorizontalAlignment = -4108
$ws.Range("O18").NumberFormat = "0"
$ws.Range("P18").Value = 8
$ws.Range("P18").HorizontalAlignment = -4108
$ws.Range("P18").NumberFormat = "0"

# ---- Selection ----
$ws.Range("O13").Select()
Write-Host "done"